$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.141153333333333
$ws.Range("H2").Value = 12.42346
$ws.Range("I2").Value = 0.2530231305454066
$ws.Range("J2").Value = 0.2530231305454066
$ws.Range("M2").Value = 0.09324
$ws.Range("O2").Value = 0.06574412111659711
$ws.Range("P2").Value = 0.06574412111659712
$ws.Range("Q2").Value = 0.3861211367999999
$ws.Range("R2").Value = 3.4750902312
$ws.Range("S2").Value = 0.01663478333987778
$ws.Range("T2").Value = 0.01663478333987778
$ws.Range("G3").Value = 4.141153333333333
$ws.Range("H3").Value = 12.42346
$ws.Range("I3").Value = 0.2530231305454066
$ws.Range("J3").Value = 0.2530231305454066
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.126335
$ws.Range("N3").Value = 0.379005
$ws.Range("O3").Value = 0.08907961755968785
$ws.Range("P3").Value = 0.08907961755968785
$ws.Range("Q3").Value = 0.5231726063666666
$ws.Range("R3").Value = 4.708553457299999
$ws.Range("S3").Value = 0.0225392037027398
$ws.Range("T3").Value = 0.0225392037027398
$ws.Range("G4").Value = 4.141153333333333
$ws.Range("H4").Value = 12.42346
$ws.Range("I4").Value = 0.2530231305454066
$ws.Range("J4").Value = 0.2530231305454066
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 1.198650666666667
$ws.Range("N4").Value = 3.595952
$ws.Range("O4").Value = 0.845176261323715
$ws.Range("P4").Value = 0.8451762613237152
$ws.Range("Q4").Value = 4.963796203768887
$ws.Range("R4").Value = 44.67416583391999
$ws.Range("S4").Value = 0.2138491435027891
$ws.Range("T4").Value = 0.2138491435027891
$ws.Range("I5").Value = 0.3583796455306321
$ws.Range("J5").Value = 0.358379645530632
$ws.Range("M5").Value = 0.09324
$ws.Range("O5").Value = 0.06574412111659711
$ws.Range("P5").Value = 0.06574412111659712
$ws.Range("Q5").Value = 0.5468984430000001
$ws.Range("R5").Value = 4.922085987000001
$ws.Range("S5").Value = 0.02356135482148902
$ws.Range("T5").Value = 0.02356135482148902
$ws.Range("I6").Value = 0.3583796455306321
$ws.Range("J6").Value = 0.358379645530632
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.126335
$ws.Range("N6").Value = 0.379005
$ws.Range("O6").Value = 0.08907961755968785
$ws.Range("P6").Value = 0.08907961755968785
$ws.Range("Q6").Value = 0.7410168897083335
$ws.Range("R6").Value = 6.669152007375001
$ws.Range("S6").Value = 0.0319243217650452
$ws.Range("T6").Value = 0.03192432176504519
$ws.Range("I7").Value = 0.3583796455306321
$ws.Range("J7").Value = 0.358379645530632
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 1.198650666666667
$ws.Range("N7").Value = 3.595952
$ws.Range("O7").Value = 0.845176261323715
$ws.Range("P7").Value = 0.8451762613237152
$ws.Range("Q7").Value = 7.030675496577778
$ws.Range("R7").Value = 63.2760794692
$ws.Range("S7").Value = 0.3028939689440979
$ws.Range("T7").Value = 0.3028939689440979
$ws.Range("G8").Value = 0.467591
$ws.Range("H8").Value = 1.402773
$ws.Range("I8").Value = 0.02856965900840602
$ws.Range("J8").Value = 0.02856965900840601
$ws.Range("M8").Value = 0.09324
$ws.Range("O8").Value = 0.06574412111659711
$ws.Range("P8").Value = 0.06574412111659712
$ws.Range("Q8").Value = 0.04359818484000001
$ws.Range("R8").Value = 0.39238366356
$ws.Range("S8").Value = 0.001878287122108525
$ws.Range("T8").Value = 0.001878287122108525
$ws.Range("G9").Value = 0.467591
$ws.Range("H9").Value = 1.402773
$ws.Range("I9").Value = 0.02856965900840602
$ws.Range("J9").Value = 0.02856965900840601
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.126335
$ws.Range("N9").Value = 0.379005
$ws.Range("O9").Value = 0.08907961755968785
$ws.Range("P9").Value = 0.08907961755968785
$ws.Range("Q9").Value = 0.059073108985
$ws.Range("R9").Value = 0.531657980865
$ws.Range("S9").Value = 0.002544974298279499
$ws.Range("T9").Value = 0.002544974298279499
$ws.Range("G10").Value = 0.467591
$ws.Range("H10").Value = 1.402773
$ws.Range("I10").Value = 0.02856965900840602
$ws.Range("J10").Value = 0.02856965900840601
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 1.198650666666667
$ws.Range("N10").Value = 3.595952
$ws.Range("O10").Value = 0.845176261323715
$ws.Range("P10").Value = 0.8451762613237152
$ws.Range("Q10").Value = 0.5604782638773333
$ws.Range("R10").Value = 5.044304374896
$ws.Range("S10").Value = 0.024146397588018
$ws.Range("T10").Value = 0.02414639758801799
$ws.Range("G11").Value = 5.892462666666667
$ws.Range("H11").Value = 17.677388
$ws.Range("I11").Value = 0.3600275649155554
$ws.Range("J11").Value = 0.3600275649155554
$ws.Range("M11").Value = 0.09324
$ws.Range("O11").Value = 0.06574412111659711
$ws.Range("P11").Value = 0.06574412111659712
$ws.Range("Q11").Value = 0.5494132190400001
$ws.Range("R11").Value = 4.94471897136
$ws.Range("S11").Value = 0.02366969583312181
$ws.Range("T11").Value = 0.02366969583312181
$ws.Range("G12").Value = 5.892462666666667
$ws.Range("H12").Value = 17.677388
$ws.Range("I12").Value = 0.3600275649155554
$ws.Range("J12").Value = 0.3600275649155554
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.126335
$ws.Range("N12").Value = 0.379005
$ws.Range("O12").Value = 0.08907961755968785
$ws.Range("P12").Value = 0.08907961755968785
$ws.Range("Q12").Value = 0.7444242709933334
$ws.Range("R12").Value = 6.69981843894
$ws.Range("S12").Value = 0.03207111779362337
$ws.Range("T12").Value = 0.03207111779362336
$ws.Range("G13").Value = 5.892462666666667
$ws.Range("H13").Value = 17.677388
$ws.Range("I13").Value = 0.3600275649155554
$ws.Range("J13").Value = 0.3600275649155554
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 1.198650666666667
$ws.Range("N13").Value = 3.595952
$ws.Range("O13").Value = 0.845176261323715
$ws.Range("P13").Value = 0.8451762613237152
$ws.Range("Q13").Value = 7.063004303708444
$ws.Range("R13").Value = 63.56703873337599
$ws.Range("S13").Value = 0.3042867512888103
$ws.Range("T13").Value = 0.3042867512888103
